$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2,3) {
    $ws.Cells.Item($r, 4).Value = -0.377                 # D
    $ws.Cells.Item($r, 7).Value = -0.7633928571428571    # G
    $ws.Cells.Item($r, 8).Value = -0.7633928571428571    # H
    $ws.Cells.Item($r, 9).Value = -2.004464285714286     # I
    $ws.Cells.Item($r, 10).Value = -1.002232142857143    # J
    $ws.Cells.Item($r, 11).Value = -23.5                 # K
    $ws.Cells.Item($r, 12).Value = -10.49107142857143    # L
    $ws.Cells.Item($r, 13).Value = 0.001                 # M
    $ws.Cells.Item($r, 14).Value = [double]"3.883495145631068e-06"   # N
    $ws.Cells.Item($r, 15).Value = [double]"-4.25531914893617e-05"   # O
    $ws.Cells.Item($r, 16).Value = 0.001                 # P
    $ws.Cells.Item($r, 17).Value = [double]"3.883495145631068e-06"   # Q
    $ws.Cells.Item($r, 18).Value = [double]"-4.25531914893617e-05"   # R
    $ws.Cells.Item($r, 19).Value = 0                     # S
    $ws.Cells.Item($r, 20).Value = 0                     # T
    $ws.Cells.Item($r, 21).Value = 0.515                 # U
    $ws.Cells.Item($r, 22).Value = 0.002                 # V
    $ws.Cells.Item($r, 23).Value = -0.8576642335766423   # W
    $ws.Cells.Item($r, 24).Value = 0.03855815202421761   # X
    $ws.Cells.Item($r, 25).Value = -0.89622238560086     # Y
    $ws.Cells.Item($r, 26).Value = 0.08291382884216762   # Z
    $ws.Cells.Item($r, 27).Value = -0.08309890435297602  # AA
    $ws.Cells.Item($r, 28).Value = 0.03855815202421761   # AB
    $ws.Cells.Item($r, 29).Value = -0.1216570563771936   # AC
    $ws.Cells.Item($r, 33).Value = -0.515                # AG
    $ws.Cells.Item($r, 36).Value = -0.002004008016032064 # AJ
    $ws.Cells.Item($r, 37).Value = -0.1539611360239163   # AK
    $ws.Cells.Item($r, 38).Value = 0                     # AL
    $ws.Cells.Item($r, 39).Value = 0                     # AM
    $ws.Cells.Item($r, 40).Value = -0                    # AN
    $ws.Cells.Item($r, 42).Value = 0.1458923512747876    # AP

    # AO and AQ are removed entirely in the new data
    $ws.Cells.Item($r, 41).ClearContents()                # AO
    $ws.Cells.Item($r, 43).ClearContents()                # AQ
}
